# Daily Scrum update - "User Case Map verbessert"
# Mark the "GUI" ToDo item as complete for Kumnig (N:Q group) and Rajic (R:U group)
# on 2018-06-29 and 2018-06-30, and add the "User Storys überarbeiten" ToDo for
# Kraschl (J:M group) plus the matching GUI ToDo entries on 2018-07-01.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (2018-06-29): Kumnig & Rajic groups get a "GUI" ToDo entry
$ws.Range("Q16").Value = "GUI"
$ws.Range("U16").Value = "GUI"

# Row 17 (2018-06-30): same "GUI" ToDo entry for Kumnig & Rajic
$ws.Range("Q17").Value = "GUI"
$ws.Range("U17").Value = "GUI"

# Row 18 (2018-07-01): Kraschl gets "User Storys überarbeiten" ToDo,
# Kumnig & Rajic again get the "GUI" ToDo entry
$ws.Range("M18").Value = "User Storys überarbeiten"
$ws.Range("Q18").Value = "GUI"
$ws.Range("U18").Value = "GUI"

# Leave the selection where the author finished editing
$ws.Range("H28").Select() | Out-Null
